# Upload new version with timestamp
# This script reproduces the change described by the diff:
#  - A new shortage-item row ("صوفي طويل جدا") is inserted into the sale sheet
#    right above the "كالونا" row (i.e. physically at row 48), shifting the
#    existing rows (كالونا, كريم فيبكس الازرق, the totals row and the footer
#    row) down by one row each.
#  - The grand-total cell is increased by the new item's price (50.0000).
#  - The footer timestamp string is bumped from 12:08 PM to 12:18 PM.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank row at row 48 - this pushes the existing rows
#    48 (كالونا), 49 (كريم فيبكس الازرق), 50 (totals) and 51 (footer)
#    down to 49, 50, 51 and 52 respectively. Excel/iron_native already
#    shifts the existing merged-cell ranges along with the rows.
$ws.Rows("48:48").Insert()

# 2) Copy the per-column formatting (styles) of the row that is now 49
#    (originally row 48 - "كالونا") into the freshly inserted row 48 so
#    that the new row visually matches the other item rows in the table.
$ws.Range("A49:Q49").Copy()
$ws.Range("A48:Q48").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Give the new row the same height/row-attributes as the other item rows.
$ws.Rows("48:48").RowHeight = 24.75

# 4) Fill in the new item's data.
$ws.Range("A48").Value = 42
$ws.Range("B48").Value = ""
$ws.Range("C48").Value = "صوفي طويل جدا"
$ws.Range("D48").Value = ""
$ws.Range("E48").Value = ""
$ws.Range("F48").Value = ""
$ws.Range("G48").Value = ""
$ws.Range("H48").Value = "5:0"
$ws.Range("I48").Value = ""
$ws.Range("J48").Value = ""
$ws.Range("K48").Value = ""
$ws.Range("L48").Value = "0"
$ws.Range("M48").Value = ""
$ws.Range("N48").Value = "50.00"
$ws.Range("O48").Value = ""
$ws.Range("P48").Value = "50.0000"
$ws.Range("Q48").Value = "1:0"

# 5) Re-create the merges for the new row (same pattern used by every
#    other item row in the table).
$ws.Range("A48:B48").Merge()
$ws.Range("C48:G48").Merge()
$ws.Range("H48:K48").Merge()
$ws.Range("L48:M48").Merge()
$ws.Range("N48:O48").Merge()

# 6) The row that used to hold the grand total (row 50) is now row 51;
#    bump its value up by the new item's price and restore the row
#    height recorded in the new workbook.
$ws.Range("P51").Value = 1892.7750000000001
$ws.Rows("51:51").RowHeight = 25.5

# 7) The footer row (used to be row 51) is now row 52; refresh the
#    timestamp text to reflect the new save time.
$ws.Range("A52").Value = "Monday, 21 July, 2025 12:18 PM"

Write-Output "Edit applied"
